# Update "想去人数" (F column) figures across all sheets to match the
# freshly regenerated gh-pages data dump (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 1656
$ws.Range("F11").Value = 1394
$ws.Range("F12").Value = 1219
$ws.Range("F14").Value = 1693
$ws.Range("F21").Value = 1151
$ws.Range("F22").Value = 1652
$ws.Range("F23").Value = 1652
$ws.Range("F30").Value = 4280
$ws.Range("F31").Value = 36
$ws.Range("F35").Value = 220
$ws.Range("F36").Value = 314
$ws.Range("F37").Value = 69

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 525
$ws.Range("F7").Value = 148041
$ws.Range("F8").Value = 148041
$ws.Range("F23").Value = 946
$ws.Range("F37").Value = 116
$ws.Range("F39").Value = 0
$ws.Range("F44").Value = 2

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1314
$ws.Range("F11").Value = 2314

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1314
$ws.Range("F9").Value = 2314
$ws.Range("F12").Value = 525
$ws.Range("F13").Value = 1938
$ws.Range("F14").Value = 148041
$ws.Range("F15").Value = 1656
$ws.Range("F18").Value = 1394
$ws.Range("F19").Value = 1219
$ws.Range("F21").Value = 1693
$ws.Range("F26").Value = 1151
$ws.Range("F27").Value = 1652
$ws.Range("F28").Value = 1652
$ws.Range("F35").Value = 4280
$ws.Range("F40").Value = 220
$ws.Range("F43").Value = 314
$ws.Range("F48").Value = 0
